$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.059.16'
$ws.Range('E2').Value = '  -2.83%  '
$ws.Range('D3').Value = '1.799.30'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.54'
$ws.Range('E5').Value = '  -3.06%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4199'
$ws.Range('E7').Value = '  -2.82%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3587'
$ws.Range('E8').Value = '  -3.50%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07110'
$ws.Range('E9').Value = '  -3.85%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8460'
$ws.Range('E10').Value = '  -4.21%  '
$ws.Range('E11').Value = '  -4.93%  '
$ws.Range('D12').Value = '1.801.93'
$ws.Range('E12').Value = '  -3.98%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.294'
$ws.Range('E13').Value = '  -4.03%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.363'
$ws.Range('E14').Value = '  -4.16%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06754'
$ws.Range('E15').Value = '  -3.04%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '80.27'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008701'
$ws.Range('E18').Value = '  -4.61%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.04'
$ws.Range('E20').Value = '  -3.84%  '
$ws.Range('D21').Value = '26.807.43'
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.054'
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.01'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = '1.963.16'
$ws.Range('E24').Value = '  -6.43%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.926'
$ws.Range('E25').Value = '  -2.93%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '152.77'
$ws.Range('E26').Value = '  -1.41%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.10'
$ws.Range('E27').Value = '  -5.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.023'
$ws.Range('E28').Value = '  -6.49%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '112.92'
$ws.Range('E29').Value = '  -2.98%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.643'
$ws.Range('E30').Value = '  -12.79%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09018'
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7238'
$ws.Range('E32').Value = '  -8.74%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.861'
$ws.Range('E33').Value = '  -4.34%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.314'
$ws.Range('E34').Value = '  -7.27%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.087'
$ws.Range('E35').Value = '  -8.36%  '
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.079'
$ws.Range('E37').Value = '  -3.18%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01906'
$ws.Range('E38').Value = '  -3.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05126'
$ws.Range('E39').Value = '  -6.54%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.1627'
$ws.Range('E40').Value = '  -4.25%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4964'
$ws.Range('E41').Value = '  -4.88%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.594'
$ws.Range('E42').Value = '  -9.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.049'
$ws.Range('E43').Value = '  -7.75%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.933'
$ws.Range('E44').Value = '  -12.82%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '104.92'
$ws.Range('E45').Value = '  -2.19%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('E47').Value = '  -4.76%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.06294'
$ws.Range('E48').Value = '  -4.37%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.4534'
$ws.Range('E49').Value = '  -6.06%  '
$ws.Range('E50').Value = '  -4.59%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.710'
$ws.Range('E51').Value = '  -8.33%  '
